# Export.xlsx template maintenance:
#  - rename the "04 Fill data to the right" sheet to "04 Horizontal fill"
#    (the four product2.* defined names that point at it follow automatically)
#  - tidy up the saved selection on that sheet
#  - leave "01 Basic Object" as the active/selected tab instead of
#    "03 Multi Table"

$wb = $excel.ActiveWorkbook

# Rename sheet 4; any defined names referencing it update in lock-step.
$ws4 = $wb.Worksheets.Item("04 Fill data to the right")
$ws4.Name = "04 Horizontal fill"

# Fix up the remembered selection on sheet 4 (without leaving it as the
# selected tab -- that honour goes to "01 Basic Object" below).
$ws4.Range("G4").Select()

# Make "01 Basic Object" the active sheet/tab (was "03 Multi Table").
$ws1 = $wb.Worksheets.Item("01 Basic Object")
$ws1.Activate()
